$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '65.402.23'
Set-TextValue 'E2' '  +4.61%  '
Set-TextValue 'D3' '3.496.63'
Set-TextValue 'E3' '  +4.22%  '
Set-TextValue 'D5' '579.77'
Set-TextValue 'E5' '  +4.18%  '
Set-TextValue 'D6' '161.70'
Set-TextValue 'E6' '  +5.57%  '
Set-TextValue 'B7' 'XRP'
Set-TextValue 'C7' 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue 'D7' '0.613'
Set-TextValue 'E7' '  +14.73%  '
Set-TextValue 'B8' 'USDC'
Set-TextValue 'C8' 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue 'D8' '0.999'
Set-TextValue 'E8' '  -0.08%  '
Set-TextValue 'D9' '3.499.15'
Set-TextValue 'E9' '  +4.24%  '
Set-TextValue 'D10' '7.28'
Set-TextValue 'E10' '  -0.64%  '
Set-TextValue 'D11' '0.126'
Set-TextValue 'E11' '  +5.44%  '
Set-TextValue 'D12' '0.447'
Set-TextValue 'E12' '  +4.49%  '
Set-TextValue 'D13' '4.093.41'
Set-TextValue 'E13' '  +3.96%  '
Set-TextValue 'E14' '  +1.31%  '
Set-TextValue 'D15' '0.0000194'
Set-TextValue 'E15' '  +5.00%  '
Set-TextValue 'D16' '28.83'
Set-TextValue 'E16' '  +8.43%  '
Set-TextValue 'D17' '65.390.62'
Set-TextValue 'E17' '  +4.34%  '
Set-TextValue 'D18' '3.500.98'
Set-TextValue 'E18' '  +3.74%  '
Set-TextValue 'D19' '6.47'
Set-TextValue 'E19' '  +5.36%  '
Set-TextValue 'D20' '14.39'
Set-TextValue 'E20' '  +3.90%  '
Set-TextValue 'D21' '384.55'
Set-TextValue 'E21' '  +3.71%  '
Set-TextValue 'D22' '8.24'
Set-TextValue 'E22' '  +4.25%  '
Set-TextValue 'D23' '0.552'
Set-TextValue 'E23' '  +6.04%  '
Set-TextValue 'D24' '72.96'
Set-TextValue 'E24' '  +3.40%  '
Set-TextValue 'E25' '  +0.72%  '
Set-TextValue 'D26' '0.0000120'
Set-TextValue 'E26' '  +6.95%  '
Set-TextValue 'D27' '10.10'
Set-TextValue 'E27' '  +8.26%  '
Set-TextValue 'D28' '0.180'
Set-TextValue 'E28' '  +3.06%  '
Set-TextValue 'E29' '  -0.04%  '
Set-TextValue 'D30' '1.53'
Set-TextValue 'E30' '  +15.93%  '
Set-TextValue 'D31' '6.26'
Set-TextValue 'E31' '  +5.02%  '
Set-TextValue 'D32' '2.06'
Set-TextValue 'E32' '  +5.16%  '
Set-TextValue 'D33' '23.69'
Set-TextValue 'E33' '  +3.68%  '
Set-TextValue 'D34' '7.28'
Set-TextValue 'E34' '  +9.68%  '
Set-TextValue 'D35' '1.59'
Set-TextValue 'E35' '  +11.55%  '
Set-TextValue 'D36' '161.73'
Set-TextValue 'E36' '  +1.99%  '
Set-TextValue 'D37' '1.93'
Set-TextValue 'E37' '  +8.13%  '
Set-TextValue 'D38' '3.024.47'
Set-TextValue 'E38' '  +4.92%  '
Set-TextValue 'D39' '0.0779'
Set-TextValue 'E39' '  +2.95%  '
Set-TextValue 'D40' '27.12'
Set-TextValue 'E40' '  +1.44%  '
Set-TextValue 'D41' '0.0324'
Set-TextValue 'E41' '  +3.65%  '
Set-TextValue 'D42' '4.59'
Set-TextValue 'E42' '  +7.98%  '
Set-TextValue 'D43' '42.84'
Set-TextValue 'E43' '  +4.85%  '
Set-TextValue 'D44' '6.50'
Set-TextValue 'E44' '  +0.35%  '
Set-TextValue 'D45' '0.782'
Set-TextValue 'E45' '  +6.00%  '
Set-TextValue 'D46' '25.99'
Set-TextValue 'E46' '  +14.64%  '
Set-TextValue 'D47' '1.12'
Set-TextValue 'E47' '  +7.11%  '
Set-TextValue 'D48' '321.70'
Set-TextValue 'E48' '  +14.12%  '
Set-TextValue 'D49' '6.78'
Set-TextValue 'E49' '  +8.03%  '
Set-TextValue 'D50' '0.110'
Set-TextValue 'E50' '  +8.61%  '
Set-TextValue 'D51' '2.21'
Set-TextValue 'E51' '  +6.78%  '
